$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1078
$ws1.Range("F4").Value = 340
$ws1.Range("F6").Value = 501
$ws1.Range("F7").Value = 8819
$ws1.Range("F8").Value = 229
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 70
$ws1.Range("F12").Value = 0

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F3").Value = 18
$ws2.Range("F5").Value = 5
$ws2.Range("F6").Value = 1

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 340
$ws4.Range("F7").Value = 0
$ws4.Range("F9").Value = 6
$ws4.Range("F10").Value = 8819
$ws4.Range("F11").Value = 229
$ws4.Range("F14").Value = 0
$ws4.Range("F16").Value = 562
$ws4.Range("F17").Value = 0
